$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 4-18 hold the "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 629
$ws1.Range("F5").Value = 171
$ws1.Range("F6").Value = 9458
$ws1.Range("F7").Value = 850
$ws1.Range("F10").Value = 1169
$ws1.Range("F11").Value = 151
$ws1.Range("F12").Value = 99
$ws1.Range("F18").Value = 1295

# Sheet "全部类型" (sheet4): same records, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 629
$ws4.Range("F6").Value = 171
$ws4.Range("F7").Value = 9458
$ws4.Range("F8").Value = 850
$ws4.Range("F11").Value = 1169
$ws4.Range("F12").Value = 151
$ws4.Range("F13").Value = 99
$ws4.Range("F19").Value = 1295
